$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "max" column (C). This shifts the old "prediction" (D) and
# "rejection-f" (E) columns left into C and D respectively.
$ws.Columns("C").Delete()

# Update the remaining data row: B2 now holds the actual prediction value.
$ws.Range("B2").Value = 319516.2110549332
